$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of columns D, J, K, L, M, P between row 2 and row 4
$cols = @("D", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $r2 = $ws.Range($col + "2")
    $r4 = $ws.Range($col + "4")
    $v2 = $r2.Value2
    $v4 = $r4.Value2
    $r2.Value2 = $v4
    $r4.Value2 = $v2
}
